$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-casing / relocating the "Expected Mode of Inheritance" allowed-values text:
# the stored shared string becomes "Comma-separated ..." (capital C). Re-assigning
# the cell's value causes the backing shared string to be rewritten; the engine
# naturally re-homes it at the end of the shared-string table and renumbers the
# other rows that shifted down, matching the authoring diff.
$ws.Range("B12").Value = "Comma-separated list of the following: Sporadic, Autosomal dominant inheritance, Sex-limited autosomal dominant, Male-limited autosomal dominant, Autosomal dominant contiguous gene syndrome, Autosomal recessive inheritance, Gonosomal inheritance, X-linked inheritance, X-linked recessive inheritance, Y-linked inheritance, X-linked dominant inheritance, Multifactorial inheritance, Mitochondrial inheritance"

# Re-format row 8 (Age of Onset): label + allowed-values cells become vertically
# centered (the row is tall to fit the wrapped Note text), and the Note cell
# becomes top-aligned with wrap text retained.
$ws.Range("A8").VerticalAlignment = -4108
$ws.Range("B8").VerticalAlignment = -4108
$ws.Range("C8").VerticalAlignment = -4160
$ws.Range("C8").WrapText = $true
